# Refresh the cryptos price/volume table (GitHub Actions daily data pull).
# Price (D) and Volume(1h) (E) columns hold plain text in the workbook, so
# numeric-looking prices are written with a leading apostrophe to keep
# Excel from auto-converting them to numbers (matches the source data,
# which stores these as text, e.g. "96.85" rather than the number 96.85).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.859.81'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '2.533.89'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''318.31'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = '''96.85'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '''0.537'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '''35.92'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').Value = '''0.0819'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '''7.56'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('E13').Value = '  -3.72%  '
$ws.Range('D14').Value = '2.921.69'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '2.499.55'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '''15.13'
$ws.Range('E16').Value = '  -4.05%  '
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '42.910.41'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '''6.87'
$ws.Range('E19').Value = '  +2.82%  '
$ws.Range('D20').Value = '''12.72'
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('D21').Value = '0.0₃0967'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '''69.70'
$ws.Range('E22').Value = '  -2.16%  '
$ws.Range('D23').Value = '''253.68'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '''26.43'
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('E28').Value = '  +2.49%  '
$ws.Range('D29').Value = '''41.18'
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('D30').Value = '''10.57'
$ws.Range('D31').Value = '''5.91'
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').Value = '''157.53'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = '''2.16'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('E34').Value = '  -2.97%  '
$ws.Range('D35').Value = '''2.71'
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('D36').Value = '''3.33'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('E39').Value = '  +7.59%  '
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('D41').Value = '''21.94'
$ws.Range('E41').Value = '  -12.56%  '
$ws.Range('D42').Value = '''0.0306'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').Value = '''3.29'
$ws.Range('E45').Value = '  -3.45%  '
$ws.Range('D46').Value = '2.009.52'
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').Value = '''9.18'
$ws.Range('E47').Value = '  +2.81%  '
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').Value = '''106.96'
$ws.Range('E49').Value = '  +4.35%  '
$ws.Range('D50').Value = '''74.91'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '2.776.96'
$ws.Range('E51').Value = '  -0.28%  '
